$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F5").Value = 3
$ws.Range("H5").Value = 3
$ws.Range("F15").Value = 78
$ws.Range("H15").Value = 78
$ws.Range("F17").Value = 46
$ws.Range("H17").Value = 46
$ws.Range("F18").Value = 44
$ws.Range("H18").Value = 44
$ws.Range("E25").Value = 19
$ws.Range("E36").Value = 91
$ws.Range("F36").Value = 40
$ws.Range("H36").Value = 40
$ws.Range("E37").Value = 46
$ws.Range("F37").Value = 25
$ws.Range("H37").Value = 25
$ws.Range("F38").Value = 15
$ws.Range("H38").Value = 15
$ws.Range("F41").Value = 15
$ws.Range("H41").Value = 15
$ws.Range("F42").Value = 14
$ws.Range("H42").Value = 14
$ws.Range("F43").Value = 14
$ws.Range("H43").Value = 14
$ws.Range("F44").Value = 12
$ws.Range("H44").Value = 12
$ws.Range("E47").Value = 55
$ws.Range("F49").Value = 34
$ws.Range("H49").Value = 34
$ws.Range("F50").Value = 7
$ws.Range("H50").Value = 7
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 8
$ws.Range("H51").Value = 8
$ws.Range("F54").Value = 1
$ws.Range("H54").Value = 1
$ws.Range("F61").Value = 9
$ws.Range("H61").Value = 9
$ws.Range("E64").Value = 33
$ws.Range("F64").Value = 17
$ws.Range("H64").Value = 17
$ws.Range("E70").Value = 40
$ws.Range("E72").Value = 36
$ws.Range("F72").Value = 17
$ws.Range("H72").Value = 17
$ws.Range("F77").Value = 18
$ws.Range("H77").Value = 18
$ws.Range("F84").Value = 3
$ws.Range("H84").Value = 3
$ws.Range("E89").Value = 32
$ws.Range("F89").Value = 13
$ws.Range("H89").Value = 13
